# Fix "Tablet Use by clinical staff" -> "Tablet use by clinical staff"
# (lowercase the "U" in "Use") and mark the corrected word with Word's
# grammar-checker proofing markers (<w:proofErr w:type="gramStart"/> ...
# <w:proofErr w:type="gramEnd"/>), same as what Word inserts when the
# grammar checker flags/revisits text after an edit.
#
# We only touch the "Use by clinical staff" portion of the run-split text
# so the preceding "Tablet " run is left completely untouched.

$d = $word.ActiveDocument

$search = "Tablet Use by clinical staff"
$text = $d.Content.Text
$idx = $text.IndexOf($search)
if ($idx -lt 0) {
    throw "Could not locate target paragraph text '$search'"
}

# Range covering "Use by clinical staff" (everything after "Tablet ")
$rng = $d.Range($idx + 7, $idx + $search.Length)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>use</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> by clinical staff</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$rng.InsertXML($xml)
